$wb = $excel.ActiveWorkbook

# The "Poland" worksheet is the template for each per-market test data
# sheet (same layout/styles); duplicate it and place the copy at the very
# end of the workbook, then turn it into the new "UK" sheet.
$src = $wb.Worksheets.Item("Poland")
$src.Copy($null, $src)

$newWs = $wb.Worksheets.Item($wb.Worksheets.Count)
$newWs.Name = "UK"

# Fill in the market-specific values (written in the same order the
# original author typed them in, so new shared strings line up).
$newWs.Range("B4").Value = "NGC-2741/T3366"
$newWs.Range("B2").Value = "UK Market"

# Leave the new sheet active/selected on cell B4, as the author had it.
$newWs.Activate()
$newWs.Range("B4").Select()
